# Added two new Mac-Addresses (10 new child rows covering two machine_ids)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 147; A = 10001; B = 10030; C = 3000166 },
    @{ Row = 148; A = 10001; B = 10030; C = 3000167 },
    @{ Row = 149; A = 10001; B = 10030; C = 3000168 },
    @{ Row = 150; A = 10001; B = 10030; C = 3000169 },
    @{ Row = 151; A = 10001; B = 10030; C = 3000170 },
    @{ Row = 152; A = 10001; B = 10031; C = 3000171 },
    @{ Row = 153; A = 10001; B = 10031; C = 3000172 },
    @{ Row = 154; A = 10001; B = 10031; C = 3000173 },
    @{ Row = 155; A = 10001; B = 10031; C = 3000174 },
    @{ Row = 156; A = 10001; B = 10031; C = 3000175 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
}

$ws.Range("A148").Select()
